# SSDM-55: fixed xls export types data.
# Adds a new "Multivalued" column (K) to the EXPERIMENT_TYPE export template:
#   - K4 header cell = "Multivalued" (bold header style, like the other headers)
#   - K5:K8 = "FALSE" (rendered via a custom TRUE/FALSE number format, left aligned)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell K4 -------------------------------------------------
# Start from the existing bold header style (row 4) so we reuse the same
# bold/size attributes, then tweak the font color to match the target look.
$ws.Range("A4").Copy() | Out-Null
$k4 = $ws.Range("K4")
$k4.PasteSpecial(-4122) | Out-Null
$k4.Value = "Multivalued"
$k4.Font.Color = 0

# --- Data cells K5:K8 -------------------------------------------------
# Literal text "FALSE" (not a boolean) displayed via a custom number format.
$dataRange = $ws.Range("K5:K8")
$dataRange.Value = "'FALSE"
$dataRange.NumberFormat = '"TRUE";"TRUE";"FALSE"'
$dataRange.HorizontalAlignment = -4131

# --- Selection / active cell ------------------------------------------
[void]$ws.Range("L7").Select()
